$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Crendetials")

# Extend the credentials table with a new row, carrying forward the
# formatting (borders/shading) used by the previous data row.
$ws.Range("B10:F10").Copy() | Out-Null
$ws.Range("B11:F11").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("B11").Value = 5
$ws.Range("C11").Value = "Putty"
$ws.Range("D11").Value = "132.148.72.192"
$ws.Range("E11").Value = "esc"
$ws.Range("F11").Value = "Esc@esc123"

# Column B in the rest of the table uses a plain bordered style; make sure
# the new cell matches that (row 10 was an outlier with an extra fill flag).
$ws.Range("B11").Borders(7).LineStyle = 1
$ws.Range("B11").Borders(8).LineStyle = 1
$ws.Range("B11").Borders(9).LineStyle = 1
$ws.Range("B11").Borders(10).LineStyle = 1

# Move the selection to the newly added cell, as seen after the edit.
$ws.Range("B11").Select() | Out-Null
